$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the payroll row for "Đào Đức Tuấn" (ID 1), currently sheet row 2.
# This shifts all subsequent employee rows up by one (their ID numbers stay
# attached to their data, they are not renumbered).
$ws.Rows.Item(2).Delete()

# --- Tham nien / khen thuong / ky luat adjustments on the remaining rows ---

# Row 3: Lưu Minh Tuấn - disciplinary fine (Tiền phạt) reduced 150000 -> 100000,
# net pay (Thực lĩnh) recalculated accordingly.
$ws.Cells.Item(3, 7).Value = 100000
$ws.Cells.Item(3, 8).Value = 10015384.615385

# Row 6: Bùi Minh Quang - bonus (Tiền thưởng) increased 700000 -> 1700000,
# net pay recalculated accordingly.
$ws.Cells.Item(6, 6).Value = 1700000
$ws.Cells.Item(6, 8).Value = 14700000

# Row 8: Nguyễn Văn Minh - work days (Ngày công) reduced 24 -> 5, salary
# (Tiền lương) and net pay recalculated accordingly.
$ws.Cells.Item(8, 3).Value = 5
$ws.Cells.Item(8, 4).Value = 2403846.1538462
$ws.Cells.Item(8, 8).Value = 2903846.1538462

# Row 14: Lê Thúy Liễu - bonus (Tiền thưởng) increased 0 -> 500000, net pay
# recalculated accordingly.
$ws.Cells.Item(14, 6).Value = 500000
$ws.Cells.Item(14, 8).Value = 8284615.3846154
